$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'247.78"
$ws.Range('D2').Style = 'Normal'
$ws.Range('D4').Value = "'5.538"
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Value = "'0.05617"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'6.483"
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Value = "'0.8085"
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Value = "'1.042"
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = "'0.1431"
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Value = "'0.07317"
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Value = "'0.03112"
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Value = "'0.02917"
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = "'0.09265"
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Value = "'0.001661"
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = "'3.228"
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Value = "'0.04731"
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = "'0.0005822"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '16OneONEWorstin24h'
$ws.Range('D18').Value = "'0.006395"
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Value = "'0.005074"
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Value = "'0.001052"
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Value = "'0.0001502"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = "'3.982"
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Value = "'3.379"
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Value = "'2.161"
$ws.Range('D24').Style = 'Normal'
$ws.Range('D26').Value = "'0.1254"
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Value = "'0.0003301"
$ws.Range('D27').Style = 'Normal'
$ws.Range('D40').Value = "'0.04144"
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Value = "'0.007130"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').Value = "'0.003505"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D43').Value = "'0.1041"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '42BKEXTokenBKK'
$ws.Range('D44').Value = "'0.008589"
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Value = "'0.00005639"
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Value = "'0.00000000751"
$ws.Range('D46').Style = 'Normal'
$ws.Range('D48').Value = "'0.01600"
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Value = "'0.00002103"
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Value = "'0.01011"
$ws.Range('D50').Style = 'Normal'

Write-Output "Applied all 42 cell updates to cryptos sheet."
